$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final dataset for the 17 companies (rows 2-18), columns A-R
$data = @(
    ('Benchling', 424.7240713084175, 4802.857225639665, 758.7945476302646, 4802.857225639665, 3394.633936788146, 162407.4561769746, 100, 1, 2, 1, 35.53137115064233, 100, 33.33333333333334, 50, 100, 54.87195980544963, 51.66666666666667),
    ('BenchSci', 320.501755284444, 0, 144.5337897812486, 320.501755284444, 118.5260448662222, 97291.88669457949, 11.85260448662222, 3, 3, 1, 1.334981932894136, 0, 100, 100, 100, 0.934487353025895, 100),
    ('Biovia', 772.7374508106509, 13253.20294340452, 2216.968971838151, 13253.20294340452, 12871.34647448357, 1147917.527963174, 100, 0, 2, 0, 100, 100, 0, 50, 0, 100, 15),
    ('CDD Vault', 181.4650918408482, 0, 79.60875022580143, 181.4650918408482, 510.4629857953323, 80665.83652537123, 51.04629857953323, 0, 1, 0, 0.2742557387983627, 44.46381412025128, 0, 0, 0, 13.53112325323424, 0),
    ('Certara', 227.8702476319986, 3430.179407605754, 253.4717113185624, 3430.179407605754, 158.5464336867516, 95399.6983528, 15.85464336867516, 3, 3, 0, 25.05905916054999, 4.540166908783551, 100, 100, 0, 18.90339148502006, 80),
    ('Dotmatics', 258.6032684038303, 3736.932106048628, 496.1372443656412, 3736.932106048628, 1488.152939379115, 298070.7640445081, 100, 1, 3, 0, 27.39930966922189, 100, 33.33333333333334, 100, 0, 49.17951676845531, 46.66666666666667),
    ('Genedata', 597.5133706123892, 2246.844304357644, 568.0612190600298, 2246.844304357644, 3186.841117373119, 173941.7822102108, 100, 3, 3, 1, 16.03126383433733, 100, 100, 100, 100, 41.22188468403613, 100),
    ('Ginkgo Bioworks', 452.672307692442, 0, 123.6009959741062, 452.672307692442, 929.6868115208051, 17964.32518467275, 92.96868115208051, 0, 1, 0, 2.343325801014551, 92.02322563591527, 0, 0, 0, 29.24729575148476, 0),
    ('KNIME', 224.9047803602669, 0, 61.55631603704023, 224.9047803602669, 4314.950036607717, 242701.3306774357, 100, 2, 3, 0, 0.6056619616943947, 100, 66.66666666666667, 100, 0, 30.42396337318608, 63.33333333333334),
    ('Labguru', 410.9886413393178, 0, 24.16462229915646, 410.9886413393178, 988.1982429404655, 79502.02923669916, 98.81982429404655, 0, 1, 1, 2.025316466846053, 98.66113377589889, 0, 0, 100, 31.0160616595619, 20),
    ('LabWare', 624.1144063085703, 3916.028672163949, 794.1433120173511, 3916.028672163949, 1296.178606936362, 243386.2914177091, 100, 0, 2, 0, 28.76565732347349, 100, 0, 50, 0, 50.13596012643144, 15),
    ('Revvity', 996.3470372789739, 7860.125762871064, 2243.929286862649, 7860.125762871064, 9877.749830401206, 1594409.267507096, 100, 1, 2, 0, 58.85561340090174, 100, 33.33333333333334, 50, 0, 71.19892938063121, 31.66666666666667),
    ('Sapio', 103.8174106722358, 0, 145.516509572369, 145.516509572369, 784.7065437552077, 60514.94778125466, 78.47065437552077, 1, 1, 0, 0, 75.57574390134782, 33.33333333333334, 0, 0, 22.67272317040435, 16.66666666666667),
    ('SciNote', 242.3934582613473, 0, 14.82963408137809, 242.3934582613473, 197.102284293974, 12828.62671180608, 19.7102284293974, 1, 1, 0, 0.7390850336404898, 8.914187307533846, 33.33333333333334, 0, 0, 3.191615715808497, 16.66666666666667),
    ('Scispot', 458.4049132667419, 0, 57.36552353828624, 458.4049132667419, 469.3446307320668, 77999.60246887438, 46.93446307320668, 2, 2, 1, 2.387060487553157, 39.79908695233114, 66.66666666666667, 50, 100, 13.61066842698655, 68.33333333333334),
    ('SimulationsPlus', 180.3881538111956, 0, 40.63204457826084, 180.3881538111956, 936.7278871083157, 82730.83416079753, 93.67278871083157, 1, 2, 1, 0.2660396585992447, 92.82200993878693, 33.33333333333334, 50, 100, 28.03283074265555, 51.66666666666667),
    ('tetrascience', 133.9565264987161, 0, 179.5862097130958, 179.5862097130958, 585.4080177240857, 82669.61396476562, 58.54080177240856, 3, 3, 1, 0.2599215377382669, 52.96605420258919, 100, 100, 100, 16.07176133719354, 100),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($rowNum, 1).Value = [string]$rowData[0]
    for ($c = 1; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = [double]$rowData[$c]
    }
}

# Remove the now-unused trailing rows (19-21) so the sheet dimension shrinks to R18
$ws.Range("A19:R21").Clear()

Write-Host "Updated sheet to 17 companies (rows 2-18)."
